# Apply the crypto-tracker data refresh (coin prices / 1h volume %, plus the
# Stacks <-> TrustWalletToken row swap near the bottom of the list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link / volume% cells, and Price cells whose text would stay text
# on its own (e.g. values with two thousands-separators like "44.202.28").
$plainUpdates = @{
    D2 = '44.202.28'
    E2 = '  +3.71%  '
    D3 = '2.255.13'
    E3 = '  +2.80%  '
    E4 = '  -0.12%  '
    E5 = '  +2.66%  '
    E6 = '  +8.48%  '
    E8 = '  -0.09%  '
    E9 = '  +3.42%  '
    E10 = '  +8.27%  '
    E11 = '  +1.71%  '
    E12 = '  +4.58%  '
    E13 = '  +2.15%  '
    D14 = '2.584.31'
    E14 = '  +2.41%  '
    E15 = '  +3.59%  '
    D16 = '2.233.93'
    E16 = '  +1.01%  '
    E17 = '  +2.06%  '
    D18 = '44.076.96'
    E18 = '  +3.69%  '
    E19 = '  +2.44%  '
    E20 = '  +0.75%  '
    E21 = '  +3.24%  '
    E22 = '  +8.99%  '
    E24 = '  +0.48%  '
    E25 = '  +0.17%  '
    E26 = '  +1.67%  '
    E27 = '  +8.34%  '
    E28 = '  -0.45%  '
    E29 = '  +2.18%  '
    E30 = '  -0.63%  '
    E31 = '  +2.42%  '
    E32 = '  +3.25%  '
    E33 = '  +10.33%  '
    E34 = '  +3.58%  '
    E35 = '  +8.52%  '
    E36 = '  +2.15%  '
    E37 = '  +12.08%  '
    E38 = '  +4.72%  '
    E39 = '  +8.21%  '
    E40 = '  +21.46%  '
    E41 = '  +3.85%  '
    E42 = '  +7.01%  '
    E45 = '  +2.14%  '
    E46 = '  +1.65%  '
    E47 = '  +2.08%  '
    E48 = '  -3.21%  '
    E49 = '  +2.66%  '
    B50 = 'TrustWalletToken'
    C50 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    E50 = '  +2.81%  '
    B51 = 'Stacks'
    C51 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    E51 = '  +25.60%  '
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Price cells whose new text parses as a plain number (e.g. "257.99"). Excel
# would silently coerce a bare .Value assignment into a numeric cell (losing
# the exact string, e.g. trailing zeros / float round-off), so mark the cell as
# Text first, write the literal string, then drop the explicit style again so
# the cell ends up with the same (unstyled) look as before the edit.
$textPriceUpdates = @{
    D5 = '257.99'
    D6 = '80.48'
    D7 = '0.627'
    D9 = '0.607'
    D10 = '43.58'
    D11 = '0.0933'
    D12 = '7.11'
    D13 = '0.104'
    D15 = '14.81'
    D17 = '0.797'
    D20 = '71.64'
    D21 = '6.08'
    D22 = '2.36'
    D23 = '235.47'
    D24 = '9.45'
    D26 = '10.87'
    D27 = '40.83'
    D28 = '3.37'
    D29 = '2.24'
    D31 = '173.44'
    D32 = '20.72'
    D33 = '0.0881'
    D34 = '5.35'
    D35 = '0.116'
    D37 = '0.0367'
    D38 = '4.53'
    D39 = '13.16'
    D40 = '2.91'
    D41 = '2.16'
    D42 = '62.96'
    D43 = '5.49'
    D44 = '0.206'
    D45 = '104.80'
    D46 = '8.59'
    D47 = '0.0995'
    D49 = '1.13'
    D50 = '1.16'
    D51 = '1.52'
}
foreach ($ref in $textPriceUpdates.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $textPriceUpdates.Keys) {
    $ws.Range($ref).Value = $textPriceUpdates[$ref]
}
foreach ($ref in $textPriceUpdates.Keys) {
    $ws.Range($ref).Style = "Normal"
}

Write-Output "Updated cryptos list."
